$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J2").Value = 2.2
$ws.Range("N2").Value = 19
$ws.Range("O2").Value = 1.14
$ws.Range("P2").Value = 5.5
$ws.Range("S2").Value = 1.25
$ws.Range("T2").Value = 3.75
$ws.Range("W2").Value = 11
$ws.Range("AB2").Value = 19
$ws.Range("AO2").Value = 8.5
$ws.Range("AR2").Value = 41
$ws.Range("AT2").Value = 3.75
$ws.Range("AW2").Value = 301
$ws.Range("AX2").Value = 6.5
